$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update subtitle text: October 2016 -> November 2016
$ws.Range("A2").Value = "Commercial Sector by Census Division and State, Year-to-Date through November 2016 (Continued)"

# Update Relative Standard Error cell values (EPM_2016_11 refresh)
$ws.Range("E4").Value = 99
$ws.Range("F4").Value = 28
$ws.Range("H4").Value = 44
$ws.Range("I4").Value = 23
$ws.Range("I5").Value = 60
$ws.Range("F6").Value = 38
$ws.Range("H6").Value = 44
$ws.Range("E7").Value = 99
$ws.Range("F7").Value = 50
$ws.Range("I7").Value = 25
$ws.Range("F8").Value = 66
$ws.Range("F9").Value = 111
$ws.Range("I9").Value = 137
$ws.Range("F10").Value = 224
$ws.Range("I10").Value = 333
$ws.Range("E11").Value = 21
$ws.Range("H11").Value = 17
$ws.Range("I11").Value = 15
$ws.Range("E12").Value = 22
$ws.Range("F12").Value = 10
$ws.Range("I12").Value = 24
$ws.Range("E13").Value = 163
$ws.Range("F13").Value = 25
$ws.Range("H13").Value = 29
$ws.Range("I13").Value = 19
$ws.Range("E14").Value = 111
$ws.Range("I14").Value = 77
$ws.Range("E15").Value = 141
$ws.Range("F15").Value = 16
$ws.Range("I15").Value = 15
$ws.Range("F16").Value = 133
$ws.Range("I16").Value = 38
$ws.Range("F17").Value = 83
$ws.Range("H17").Value = 99
$ws.Range("I17").Value = 32
$ws.Range("F18").Value = 14
$ws.Range("I18").Value = 19
$ws.Range("E19").Value = 141
$ws.Range("F19").Value = 119
$ws.Range("I19").Value = 103
$ws.Range("F20").Value = 63
$ws.Range("I20").Value = 75
$ws.Range("E21").Value = 204
$ws.Range("F21").Value = 31
$ws.Range("H21").Value = 72
$ws.Range("I21").Value = 25
$ws.Range("I22").Value = 53
$ws.Range("F23").Value = 38
$ws.Range("H23").Value = 72
$ws.Range("I23").Value = 51
$ws.Range("E24").Value = 204
$ws.Range("F24").Value = 204
$ws.Range("I25").Value = 108
$ws.Range("I26").Value = 900
$ws.Range("I27").Value = 758
$ws.Range("E28").Value = 24
$ws.Range("H28").Value = 16
$ws.Range("I28").Value = 20
$ws.Range("E29").Value = 199
$ws.Range("F29").Value = 113
$ws.Range("I29").Value = 113
$ws.Range("I30").Value = 156
$ws.Range("E31").Value = 192
$ws.Range("F31").Value = 54
$ws.Range("I31").Value = 86
$ws.Range("E32").Value = 136
$ws.Range("F32").Value = 99
$ws.Range("I32").Value = 61
$ws.Range("E33").Value = 100
$ws.Range("F33").Value = 51
$ws.Range("I33").Value = 54
$ws.Range("E34").Value = 25
$ws.Range("F34").Value = 24
$ws.Range("I34").Value = 17
$ws.Range("I35").Value = 297
$ws.Range("F36").Value = 13
$ws.Range("H36").Value = 16
$ws.Range("I36").Value = 12
$ws.Range("E37").Value = 142
$ws.Range("F37").Value = 142
$ws.Range("I37").Value = 93
$ws.Range("I38").Value = 348
$ws.Range("E39").Value = 142
$ws.Range("F39").Value = 142
$ws.Range("I39").Value = 97
$ws.Range("E40").Value = 169
$ws.Range("F40").Value = 50
$ws.Range("I40").Value = 27
$ws.Range("F41").Value = 186
$ws.Range("I41").Value = 306
$ws.Range("I42").Value = 86
$ws.Range("I43").Value = 253
$ws.Range("E44").Value = 169
$ws.Range("F44").Value = 51
$ws.Range("I44").Value = 28
$ws.Range("E45").Value = 26
$ws.Range("F45").Value = 27
$ws.Range("I45").Value = 16
$ws.Range("E46").Value = 57
$ws.Range("F46").Value = 57
$ws.Range("I46").Value = 26
$ws.Range("E47").Value = 65
$ws.Range("F47").Value = 65
$ws.Range("I47").Value = 187
$ws.Range("F48").Value = 131
$ws.Range("I48").Value = 131
$ws.Range("E49").Value = 33
$ws.Range("F49").Value = 33
$ws.Range("I49").Value = 26
$ws.Range("F50").Value = 320
$ws.Range("I50").Value = 30
$ws.Range("F51").Value = 100
$ws.Range("I51").Value = 34
$ws.Range("E52").Value = 24
$ws.Range("F52").Value = 10
$ws.Range("I52").Value = 7
$ws.Range("E53").Value = 24
$ws.Range("F53").Value = 10
$ws.Range("I53").Value = 7
$ws.Range("F54").Value = 87
$ws.Range("I54").Value = 101
$ws.Range("I55").Value = 131
$ws.Range("F56").Value = 8
$ws.Range("I56").Value = 8
$ws.Range("F57").Value = 54
$ws.Range("I57").Value = 30
$ws.Range("E59").Value = 12
$ws.Range("F59").Value = 6
